$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3310
$ws.Range("J64").Value = 3580
$ws.Range("L64").Value = 3580
$ws.Range("N64").Value = -4076
$ws.Range("H67").Value = 3310
$ws.Range("J67").Value = 3580
$ws.Range("L67").Value = 3580
$ws.Range("N67").Value = -5296
$ws.Range("H76").Value = 3149.5
$ws.Range("I76").Value = 2999.5
$ws.Range("J76").Value = 3299.5
$ws.Range("K76").Value = 2999.5
$ws.Range("L76").Value = 3299.5
$ws.Range("M76").Value = -2684.5
$ws.Range("N76").Value = -3929.5
$ws.Range("H79").Value = 3149.5
$ws.Range("I79").Value = 2999.5
$ws.Range("J79").Value = 3299.5
$ws.Range("K79").Value = 2999.5
$ws.Range("L79").Value = 3299.5
$ws.Range("M79").Value = -1907.5
$ws.Range("N79").Value = -5483.5
$ws.Range("H98").Value = 4069.1155
$ws.Range("I98").Value = 3883.7222
$ws.Range("J98").Value = 4486.25
$ws.Range("K98").Value = 3883.7222
$ws.Range("L98").Value = 4486.25
$ws.Range("M98").Value = -2385.7222
$ws.Range("N98").Value = -7482.25
$ws.Range("H122").Value = 4069.1155
$ws.Range("I122").Value = 3883.7222
$ws.Range("J122").Value = 4486.25
$ws.Range("K122").Value = 11651.1666
$ws.Range("L122").Value = 13458.75
$ws.Range("M122").Value = -9201.1666
$ws.Range("N122").Value = -18358.75
$ws.Range("H141").Value = 4002677.8
$ws.Range("I141").Value = 7000797.5
$ws.Range("J141").Value = 5185
$ws.Range("K141").Value = 21002392.5
$ws.Range("L141").Value = 15555
$ws.Range("M141").Value = -20997212.5
$ws.Range("N141").Value = -25915

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5882.6577
$ws.Range("I32").Value = 3238.0862
$ws.Range("K32").Value = 3238.0862
$ws.Range("M32").Value = -2951.0862
$ws.Range("H45").Value = 1624.1666
$ws.Range("I45").Value = 1249.5
$ws.Range("J45").Value = 1811.5
$ws.Range("K45").Value = 1249.5
$ws.Range("L45").Value = 1811.5
$ws.Range("M45").Value = -872.5
$ws.Range("N45").Value = -2565.5
$ws.Range("H63").Value = 2113
$ws.Range("I63").Value = 1726
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 1726
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -1040
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 2113
$ws.Range("I66").Value = 1726
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 8630
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -5198
$ws.Range("N66").Value = -19364
$ws.Range("H97").Value = 2744.2222
$ws.Range("I97").Value = 2774.75
$ws.Range("K97").Value = 2774.75
$ws.Range("M97").Value = -2278.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 798
$ws.Range("I94").Value = 388.2857
$ws.Range("K94").Value = 388.2857
$ws.Range("M94").Value = 62.71429999999998
$ws.Range("H105").Value = 2770.8948
$ws.Range("I105").Value = 2514.8125
$ws.Range("J105").Value = 4136.6665
$ws.Range("K105").Value = 2514.8125
$ws.Range("L105").Value = 4136.6665
$ws.Range("M105").Value = -767.8125
$ws.Range("N105").Value = -7630.6665
$ws.Range("H134").Value = 7559.375
$ws.Range("I134").Value = 10524.429
$ws.Range("J134").Value = 3408.3
$ws.Range("K134").Value = 31573.287
$ws.Range("L134").Value = 10224.9
$ws.Range("M134").Value = -29038.287
$ws.Range("N134").Value = -15294.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3346758.2
$ws.Range("I58").Value = 7249308.5
$ws.Range("J58").Value = 1715
$ws.Range("K58").Value = 7249308.5
$ws.Range("L58").Value = 1715
$ws.Range("M58").Value = -7249105.5
$ws.Range("N58").Value = -2121
$ws.Range("H105").Value = 1072.5714
$ws.Range("I105").Value = 1101.5454
$ws.Range("K105").Value = 1101.5454
$ws.Range("M105").Value = 645.4546
$ws.Range("H134").Value = 3455.2727
$ws.Range("I134").Value = 2913.625
$ws.Range("K134").Value = 8740.875
$ws.Range("M134").Value = -6205.875
$ws.Range("H136").Value = 3346758.2
$ws.Range("I136").Value = 7249308.5
$ws.Range("J136").Value = 1715
$ws.Range("K136").Value = 21747925.5
$ws.Range("L136").Value = 5145
$ws.Range("M136").Value = -21745375.5
$ws.Range("N136").Value = -10245

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1237131.4
$ws.Range("I4").Value = 1281406
$ws.Range("K4").Value = 3844218
$ws.Range("M4").Value = -3844106
$ws.Range("H10").Value = 193.8
$ws.Range("I10").Value = 192.5
$ws.Range("K10").Value = 577.5
$ws.Range("M10").Value = -438.5
$ws.Range("H129").Value = 25247.467
$ws.Range("I129").Value = 700.5714
$ws.Range("J129").Value = 32718.262
$ws.Range("K129").Value = 2101.7142
$ws.Range("L129").Value = 98154.78599999999
$ws.Range("M129").Value = 2898.2858
$ws.Range("N129").Value = -108154.786
$ws.Range("H131").Value = 16820.268
$ws.Range("I131").Value = 297.6
$ws.Range("K131").Value = 892.8000000000001
$ws.Range("M131").Value = 4147.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4355.591
$ws.Range("I70").Value = 4136.615
$ws.Range("K70").Value = 4136.615
$ws.Range("M70").Value = -3866.615
$ws.Range("H73").Value = 4355.591
$ws.Range("I73").Value = 4136.615
$ws.Range("K73").Value = 4136.615
$ws.Range("M73").Value = -3200.615
$ws.Range("H80").Value = 2787.0667
$ws.Range("I80").Value = 2676.4614
$ws.Range("K80").Value = 2676.4614
$ws.Range("M80").Value = -1678.4614
$ws.Range("H83").Value = 2787.0667
$ws.Range("I83").Value = 2676.4614
$ws.Range("K83").Value = 13382.307
$ws.Range("M83").Value = -8390.307000000001
$ws.Range("H97").Value = 1051.9259
$ws.Range("I97").Value = 1043.619
$ws.Range("J97").Value = 1081
$ws.Range("K97").Value = 1043.619
$ws.Range("L97").Value = 1081
$ws.Range("M97").Value = -547.6189999999999
$ws.Range("N97").Value = -2073
$ws.Range("H132").Value = 1834956.4
$ws.Range("I132").Value = 2265397.5
$ws.Range("K132").Value = 6796192.5
$ws.Range("M132").Value = -6793662.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9080.137000000001
$ws.Range("I40").Value = 9848.0625
$ws.Range("J40").Value = 7032.3335
$ws.Range("K40").Value = 9848.0625
$ws.Range("L40").Value = 7032.3335
$ws.Range("M40").Value = -9712.0625
$ws.Range("N40").Value = -7304.3335
$ws.Range("H82").Value = 2370
$ws.Range("I82").Value = 1745
$ws.Range("J82").Value = 2995
$ws.Range("K82").Value = 1745
$ws.Range("L82").Value = 2995
$ws.Range("M82").Value = -1384
$ws.Range("N82").Value = -3717
$ws.Range("H85").Value = 2370
$ws.Range("I85").Value = 1745
$ws.Range("J85").Value = 2995
$ws.Range("K85").Value = 1745
$ws.Range("L85").Value = 2995
$ws.Range("M85").Value = -497
$ws.Range("N85").Value = -5491
$ws.Range("H93").Value = 1173.1904
$ws.Range("I93").Value = 648.05884
$ws.Range("J93").Value = 3405
$ws.Range("K93").Value = 648.05884
$ws.Range("L93").Value = 3405
$ws.Range("M93").Value = 599.94116
$ws.Range("N93").Value = -5901

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4332.3335
$ws.Range("J62").Value = 3998.5
$ws.Range("L62").Value = 3998.5
$ws.Range("N62").Value = -5246.5
$ws.Range("H65").Value = 4332.3335
$ws.Range("J65").Value = 3998.5
$ws.Range("L65").Value = 19992.5
$ws.Range("N65").Value = -26232.5
$ws.Range("H81").Value = 4184.8335
$ws.Range("I81").Value = 1777.75
$ws.Range("J81").Value = 8999
$ws.Range("K81").Value = 3555.5
$ws.Range("L81").Value = 17998
$ws.Range("M81").Value = -2494.5
$ws.Range("N81").Value = -20120
$ws.Range("H84").Value = 4184.8335
$ws.Range("I84").Value = 1777.75
$ws.Range("J84").Value = 8999
$ws.Range("K84").Value = 17777.5
$ws.Range("L84").Value = 89990
$ws.Range("M84").Value = -12473.5
$ws.Range("N84").Value = -100598
$ws.Range("H125").Value = 129999.95
$ws.Range("J125").Value = 129999.95
$ws.Range("L125").Value = 129999.95
$ws.Range("N125").Value = -139839.95
$ws.Range("H126").Value = 1860.0834
$ws.Range("I126").Value = 1596.8823
$ws.Range("K126").Value = 4790.6469
$ws.Range("M126").Value = -2320.6469
